$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain as text, matching the source format
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '61.326.01'
$ws.Range("E2").Value = '  +7.97%  '
$ws.Range("D3").Value = '3.408.77'
$ws.Range("E3").Value = '  +5.31%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '412.78'
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").Value = '122.76'
$ws.Range("E6").Value = '  +13.70%  '
$ws.Range("D7").Value = '3.402.63'
$ws.Range("E7").Value = '  +5.27%  '
$ws.Range("D8").Value = '0.582'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("E10").Value = '  +4.07%  '
$ws.Range("E11").Value = '  +17.77%  '
$ws.Range("D12").Value = '41.59'
$ws.Range("E12").Value = '  +5.97%  '
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").Value = '3.962.08'
$ws.Range("E14").Value = '  +5.77%  '
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").Value = '19.68'
$ws.Range("E16").Value = '  +4.09%  '
$ws.Range("D17").Value = '3.421.12'
$ws.Range("E17").Value = '  +5.65%  '
$ws.Range("D18").Value = '61.377.72'
$ws.Range("E18").Value = '  +8.32%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = '10.94'
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("E21").Value = '  +6.45%  '
$ws.Range("D22").Value = '3.38'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").Value = '13.01'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '299.07'
$ws.Range("E24").Value = '  +2.69%  '
$ws.Range("D25").Value = '76.04'
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("D26").Value = '3.14'
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("D27").Value = '30.81'
$ws.Range("E27").Value = '  +9.87%  '
$ws.Range("D28").Value = '8.24'
$ws.Range("E28").Value = '  +14.23%  '
$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").Value = '43.06'
$ws.Range("E32").Value = '  +3.81%  '
$ws.Range("E33").Value = '  +5.39%  '
$ws.Range("D34").Value = '11.47'
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '2.53'
$ws.Range("E36").Value = '  +18.67%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '52.19'
$ws.Range("E38").Value = '  +1.53%  '
$ws.Range("D39").Value = '3.55'
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  +2.05%  '
$ws.Range("D42").Value = '1.97'
$ws.Range("E42").Value = '  +5.04%  '
$ws.Range("D46").Value = '3.98'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("D47").Value = '0.283'
$ws.Range("E47").Value = '  +2.58%  '
$ws.Range("D48").Value = '22.11'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("D49").Value = '2.19'
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").Value = '2.196.90'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("D51").Value = '3.757.65'
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '134.57'
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.122'
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '17.50'
$ws.Range("E45").Value = '  +3.56%  '
